$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (was "GossA-HW25.xpc")
$ws.Name = "GossA"

# Correct tiny floating point differences in row 13
$ws.Range("D13").Value = 0.9943698965328298
$ws.Range("H13").Value = 0.9943698965328298
$ws.Range("L13").Value = 0.9930084784487491
$ws.Range("N13").Value = 0.9944867442981457

# Prepare row 16 formatting by copying row 15's label-cell style (bold,
# bordered, centered) onto A16 before filling in the new data
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)

# Add new row 16 of averaged-intensity data (HKL index 14, HexGrid-60degTilt5degRes)
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 1.086662413002611
$ws.Range("D16").Value = 0.9491496470760641
$ws.Range("E16").Value = 0.9898537016774713
$ws.Range("F16").Value = 0.9725266603351137
$ws.Range("G16").Value = 1.086662413002611
$ws.Range("H16").Value = 0.9491496470760641
$ws.Range("I16").Value = 1.019375589268239
$ws.Range("J16").Value = 0.9605942201170385
$ws.Range("K16").Value = 1.018040081265707
$ws.Range("L16").Value = 0.9504567340761417
$ws.Range("M16").Value = 1.086662413002611
$ws.Range("N16").Value = 0.9695016743767677
$ws.Range("O16").Value = 0.999548105522815
$ws.Range("P16").Value = 0.9933323808522984
